# Actualización automática 2025-06-01 08:00:06
#
# Applies the monthly roll-forward to "GUERRERO FAREZ FABIAN MAURICIO":
#   - Sheet "VENTAS POR GRUPO": clears every product/group total back to 0
#     (a brand-new reporting period with no sales yet) and the "N de 50"
#     progress row is reset to "0 de 50". Column I is narrowed to match.
#   - Sheet "VENTA MENSUAL": the monthly columns roll forward one month
#     (the old "marzo" column's data becomes "febrero"'s position's data,
#     etc.) - i.e. column C takes column D's old values, D takes E's, E
#     takes F's, and the new trailing month F starts at 0. Headers shift
#     the same way, with the new month "junio" introduced in column F.

$wb = $excel.ActiveWorkbook

# Excel's ColumnWidth setter round-trips through a pixel-based internal
# width that lands 5/6 of a character above whatever gets written to the
# saved <col width="..."/> attribute, so compensate for that offset to
# land on the exact target width.
$widthOffset = 0.8333333333333334

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Column I (9) narrows from 12 to 9 characters.
$ws1.Columns.Item(9).ColumnWidth = 9 - $widthOffset

# Every non-zero product-group total in the data rows goes to 0.
$cellsToZero = @(
    "C4","D4","E4","F4","K4","L4","M4",
    "C5","H5","K5","L5",
    "C6","L6",
    "N10",
    "L13",
    "D14","L14","N14",
    "L15",
    "C16","E16","F16","K16","L16",
    "L18",
    "C24","E24","F24","J24","K24","L24",
    "L27",
    "E28","F28","L28","N28",
    "E30",
    "D31","L31",
    "C33","D33","K33","L33",
    "D34","E34","F34","K34","L34",
    "I35",
    "E38","F38","J38","K38","L38",
    "J39",
    "D40","E40","F40","L40",
    "L45",
    "C46","D46",
    "L50"
)
foreach ($ref in $cellsToZero) {
    $ws1.Range($ref).Value2 = 0
}

# Row 52 holds the "<n> de 50" counters per column; all counts drop to 0.
foreach ($col in @("C","D","E","F","G","H","I","J","K","L","M","N")) {
    $ws1.Range("$col`52").Value2 = "0 de 50"
}

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Column widths shuffle along with the shifted months.
$ws2.Columns.Item(4).ColumnWidth = 13 - $widthOffset
$ws2.Columns.Item(5).ColumnWidth = 14 - $widthOffset
$ws2.Columns.Item(6).ColumnWidth = 11 - $widthOffset

# Header months roll forward by one: marzo, abril, mayo, junio.
$ws2.Range("C1").Value2 = "marzo"
$ws2.Range("D1").Value2 = "abril"
$ws2.Range("E1").Value2 = "mayo"
$ws2.Range("F1").Value2 = "junio"

# Every data (and totals) row 2..52 shifts its monthly figures left by one
# column, with the new trailing month (F) starting empty (0).
for ($r = 2; $r -le 52; $r++) {
    $oldD = $ws2.Cells.Item($r, 4).Value2
    $oldE = $ws2.Cells.Item($r, 5).Value2
    $oldF = $ws2.Cells.Item($r, 6).Value2

    $ws2.Cells.Item($r, 3).Value2 = $oldD
    $ws2.Cells.Item($r, 4).Value2 = $oldE
    $ws2.Cells.Item($r, 5).Value2 = $oldF
    $ws2.Cells.Item($r, 6).Value2 = 0
}
